$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Workbook-level window / fileVersion bookkeeping
# ---------------------------------------------------------------------------
$wb.Windows.Item(1).Left   = 2880
$wb.Windows.Item(1).Top    = 1160

# ---------------------------------------------------------------------------
# guess_me (sheet1): selection moves from C11 to E6
# ---------------------------------------------------------------------------
$wsGuessMe = $wb.Worksheets.Item("guess_me")
$wsGuessMe.Range("E6").Select()

# ---------------------------------------------------------------------------
# guess_max (sheet2): selection sqref shrinks from A1:E1048576 to A1:E65536
# ---------------------------------------------------------------------------
$wsGuessMax = $wb.Worksheets.Item("guess_max")
$wsGuessMax.Range("A1:E65536").Select()

# ---------------------------------------------------------------------------
# numeric_coercion (sheet6): A3 becomes the text "72", right aligned, bordered
# on the left with a medium border, and the selection moves to A3.
# (processed first so the new shared string "72" lands at index 31 and the
# new cell style lands at cellXfs index 3, matching the target workbook)
# ---------------------------------------------------------------------------
$wsNumeric = $wb.Worksheets.Item("numeric_coercion")
$wsNumeric.Range("A3").ClearFormats()
$wsNumeric.Range("A3").NumberFormat = "@"
$wsNumeric.Range("A3").HorizontalAlignment = -4152
$wsNumeric.Range("A3").Borders.Item(7).Weight = -4138
$wsNumeric.Range("A3").Value = "72"
$wsNumeric.Range("A3").Select()

# ---------------------------------------------------------------------------
# date_coercion (sheet5): add the bestFit column def, and A7 changes from a
# date-formatted 39529 to a General-formatted 39448.
# ---------------------------------------------------------------------------
$wsDate = $wb.Worksheets.Item("date_coercion")
$wsDate.Columns.Item(1).ColumnWidth = 11
$wsDate.Range("A7").ClearFormats()
$wsDate.Range("A7").NumberFormat = "General"
$wsDate.Range("A7").Value = 39448

# ---------------------------------------------------------------------------
# logical_coercion (sheet4): rebuilt with an added "explanation" column (B)
# and several new example rows covering string/quote-prefixed logicals.
# ---------------------------------------------------------------------------
$wsLogical = $wb.Worksheets.Item("logical_coercion")

$wsLogical.Range("A1").Value = "logical"
$wsLogical.Range("B1").Value = "explanation"

$wsLogical.Range("A2").Value = $true
$wsLogical.Range("B2").Value = "static logical"

$wsLogical.Range("A3").Value = $false
$wsLogical.Range("B3").Value = "static logical"

$wsLogical.Range("A4").Formula = "=TRUE()"
$wsLogical.Range("B4").Value = "formula logical"

$wsLogical.Range("A5").Formula = "=FALSE()"
$wsLogical.Range("B5").Value = "formula logical"

$wsLogical.Range("A6").Formula = '="true"'
$wsLogical.Range("B6").Value = "string logical"

$wsLogical.Range("A7").ClearFormats()
$wsLogical.Range("A7").Formula = '="false"'
$wsLogical.Range("B7").Value = "string logical"

$wsLogical.Range("A8").Value = "'true"
$wsLogical.Range("B8").Value = "string logical"

$wsLogical.Range("A9").Value = "'false"
$wsLogical.Range("B9").Value = "string logical"

$wsLogical.Range("A10").Value = "T"
$wsLogical.Range("B10").Value = "string logical"

$wsLogical.Range("A11").Value = "F"
$wsLogical.Range("B11").Value = "string logical"

$wsLogical.Range("A12").Value = "'True"
$wsLogical.Range("B12").Value = "string logical"

$wsLogical.Range("A13").Value = "'False"
$wsLogical.Range("B13").Value = "string logical"

$wsLogical.Range("A14").Value = "cabbage"
$wsLogical.Range("B14").Value = "string not logical"

$wsLogical.Range("B15").Value = "blank"

$wsLogical.Range("A16").Value = 0
$wsLogical.Range("B16").Value = "numeric"

$wsLogical.Range("A17").Value = 1
$wsLogical.Range("B17").Value = "numeric"

$wsLogical.Range("A18").NumberFormat = "m/d/yyyy"
$wsLogical.Range("A18").Value = 40908
$wsLogical.Range("B18").Value = "date"

$wsLogical.Range("B18").Select()

# ---------------------------------------------------------------------------
# text_coercion (sheet7): add matching "explanation" column (B) alongside the
# existing data in column A.
# ---------------------------------------------------------------------------
$wsText = $wb.Worksheets.Item("text_coercion")

$wsText.Range("B1").Value = "explanation"
$wsText.Range("B2").Value = "text"
$wsText.Range("B3").Value = "blank"
$wsText.Range("B4").Value = "logical F"
$wsText.Range("B5").Value = "boolean"
$wsText.Range("B6").Value = "floating point"
$wsText.Range("B7").Value = "date"
$wsText.Range("B8").Value = "text"

$wsText.Range("B9").Select()
